$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 4524.8945
$ws.Range("J45").Value = 4524.8945
$ws.Range("L45").Value = 13574.6835
$ws.Range("N45").Value = -13958.6835
$ws.Range("H55").Value = 72329.71000000001
$ws.Range("J55").Value = 1920.5
$ws.Range("L55").Value = 1920.5
$ws.Range("N55").Value = -2348.5
$ws.Range("H112").Value = 1972.7273
$ws.Range("J112").Value = 1972.7273
$ws.Range("L112").Value = 5918.1819
$ws.Range("N112").Value = -8134.1819
$ws.Range("H125").Value = 12350001
$ws.Range("J125").Value = 18522852
$ws.Range("L125").Value = 166705668
$ws.Range("N125").Value = -166710588
$ws.Range("H137").Value = 2358.3928
$ws.Range("I137").Value = 2347
$ws.Range("K137").Value = 7041
$ws.Range("M137").Value = -4491

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 50012500
$ws.Range("I36").Value = 50012500
$ws.Range("K36").Value = 50012500
$ws.Range("M36").Value = -50012154
$ws.Range("H74").Value = 1520.6538
$ws.Range("I74").Value = 1501.48
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1501.48
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -627.48
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1520.6538
$ws.Range("I77").Value = 1501.48
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 7507.4
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -3139.4
$ws.Range("N77").Value = -18736
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 39866.63
$ws.Range("I134").Value = 2915.96
$ws.Range("K134").Value = 8747.880000000001
$ws.Range("M134").Value = -6212.880000000001

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41796.12
$ws.Range("I31").Value = 1300.174
$ws.Range("J31").Value = 507499.5
$ws.Range("K31").Value = 1300.174
$ws.Range("L31").Value = 507499.5
$ws.Range("M31").Value = -1005.174
$ws.Range("N31").Value = -508089.5
$ws.Range("H34").Value = 41796.12
$ws.Range("I34").Value = 1300.174
$ws.Range("J34").Value = 507499.5
$ws.Range("K34").Value = 1300.174
$ws.Range("L34").Value = 507499.5
$ws.Range("M34").Value = -1098.174
$ws.Range("N34").Value = -507903.5
$ws.Range("H94").Value = 913.5714
$ws.Range("J94").Value = 939
$ws.Range("L94").Value = 939
$ws.Range("N94").Value = -1841
$ws.Range("H99").Value = 5172.5
$ws.Range("J99").Value = 5944.5713
$ws.Range("L99").Value = 5944.5713
$ws.Range("N99").Value = -8940.5713
$ws.Range("H122").Value = 4203.95
$ws.Range("J122").Value = 4568.125
$ws.Range("L122").Value = 13704.375
$ws.Range("N122").Value = -18604.375
$ws.Range("H126").Value = 5172.5
$ws.Range("J126").Value = 5944.5713
$ws.Range("L126").Value = 17833.7139
$ws.Range("N126").Value = -22773.7139
$ws.Range("H132").Value = 2267.3333
$ws.Range("I132").Value = 1906.6666
$ws.Range("J132").Value = 3169
$ws.Range("K132").Value = 5719.9998
$ws.Range("L132").Value = 9507
$ws.Range("M132").Value = -3189.9998
$ws.Range("N132").Value = -14567
$ws.Range("H141").Value = 328498
$ws.Range("J141").Value = 342378.12
$ws.Range("L141").Value = 342378.12
$ws.Range("N141").Value = -352738.12

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 146675.14
$ws.Range("J34").Value = 171087.67
$ws.Range("L34").Value = 513263.01
$ws.Range("N34").Value = -513431.01
$ws.Range("H50").Value = 41959.082
$ws.Range("J50").Value = 71508.42999999999
$ws.Range("L50").Value = 214525.29
$ws.Range("N50").Value = -215487.29
$ws.Range("H53").Value = 41959.082
$ws.Range("J53").Value = 71508.42999999999
$ws.Range("L53").Value = 214525.29
$ws.Range("N53").Value = -215487.29
$ws.Range("H60").Value = 351.42856
$ws.Range("H70").Value = 2548.4
$ws.Range("I70").Value = 2548.4
$ws.Range("K70").Value = 7645.200000000001
$ws.Range("M70").Value = -7330.200000000001
$ws.Range("H73").Value = 2548.4
$ws.Range("I73").Value = 2548.4
$ws.Range("K73").Value = 7645.200000000001
$ws.Range("M73").Value = -6553.200000000001
$ws.Range("H121").Value = 501704.75
$ws.Range("I121").Value = 1642.4286
$ws.Range("K121").Value = 4927.2858
$ws.Range("M121").Value = -3617.2858
$ws.Range("H122").Value = 41473.52
$ws.Range("I122").Value = 825.7143
$ws.Range("J122").Value = 57281
$ws.Range("K122").Value = 7431.428699999999
$ws.Range("L122").Value = 515529
$ws.Range("M122").Value = -4981.428699999999
$ws.Range("N122").Value = -520429
$ws.Range("H131").Value = 2930.6304
$ws.Range("I131").Value = 2018.6666
$ws.Range("J131").Value = 3152.4595
$ws.Range("K131").Value = 6055.9998
$ws.Range("L131").Value = 9457.378499999999
$ws.Range("M131").Value = -1015.9998
$ws.Range("N131").Value = -19537.3785

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 504.5
$ws.Range("J17").Value = 504.5
$ws.Range("L17").Value = 504.5
$ws.Range("N17").Value = -840.5
$ws.Range("H80").Value = 743995.1
$ws.Range("I80").Value = 670571.4
$ws.Range("J80").Value = 835774.75
$ws.Range("K80").Value = 670571.4
$ws.Range("L80").Value = 835774.75
$ws.Range("M80").Value = -669573.4
$ws.Range("N80").Value = -837770.75
$ws.Range("H83").Value = 743995.1
$ws.Range("I83").Value = 670571.4
$ws.Range("J83").Value = 835774.75
$ws.Range("K83").Value = 3352857
$ws.Range("L83").Value = 4178873.75
$ws.Range("M83").Value = -3347865
$ws.Range("N83").Value = -4188857.75
$ws.Range("H102").Value = 5224
$ws.Range("I102").Value = 2558.8572
$ws.Range("K102").Value = 2558.8572
$ws.Range("M102").Value = -936.8571999999999
$ws.Range("H132").Value = 42901.19
$ws.Range("I132").Value = 4617.24
$ws.Range("K132").Value = 13851.72
$ws.Range("M132").Value = -11321.72

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6613.8125
$ws.Range("I7").Value = 6074.636
$ws.Range("J7").Value = 7800
$ws.Range("K7").Value = 6074.636
$ws.Range("L7").Value = 7800
$ws.Range("M7").Value = -5962.636
$ws.Range("N7").Value = -8024
$ws.Range("H61").Value = 2910.682
$ws.Range("I61").Value = 3001.4285
$ws.Range("K61").Value = 3001.4285
$ws.Range("M61").Value = -2799.4285
$ws.Range("H68").Value = 2774.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 2774.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 1360.6154
$ws.Range("I82").Value = 1419.4445
$ws.Range("J82").Value = 1228.25
$ws.Range("K82").Value = 1419.4445
$ws.Range("L82").Value = 1228.25
$ws.Range("M82").Value = -1058.4445
$ws.Range("N82").Value = -1950.25
$ws.Range("H85").Value = 1360.6154
$ws.Range("I85").Value = 1419.4445
$ws.Range("J85").Value = 1228.25
$ws.Range("K85").Value = 1419.4445
$ws.Range("L85").Value = 1228.25
$ws.Range("M85").Value = -171.4445000000001
$ws.Range("N85").Value = -3724.25
$ws.Range("H113").Value = 2910.682
$ws.Range("I113").Value = 3001.4285
$ws.Range("K113").Value = 3001.4285
$ws.Range("M113").Value = -831.4285
$ws.Range("H126").Value = 6613.8125
$ws.Range("I126").Value = 6074.636
$ws.Range("J126").Value = 7800
$ws.Range("K126").Value = 18223.908
$ws.Range("L126").Value = 23400
$ws.Range("M126").Value = -15753.908
$ws.Range("N126").Value = -28340

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11538.333
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 11538.333
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 11538.333
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -13410.333
$ws.Range("H77").Value = 11538.333
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 11538.333
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 34614.999
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -43974.999
$ws.Range("H100").Value = 1048.8462
$ws.Range("I100").Value = 914.3
$ws.Range("K100").Value = 1828.6
$ws.Range("M100").Value = -1287.6
$ws.Range("H127").Value = 45000
$ws.Range("J127").Value = 45000
$ws.Range("L127").Value = 45000
$ws.Range("N127").Value = -54920
$ws.Range("H130").Value = 39690
$ws.Range("J130").Value = 39690
$ws.Range("L130").Value = 39690
$ws.Range("N130").Value = -49730
$ws.Range("H132").Value = 19578.525
$ws.Range("I132").Value = 2348.9575
$ws.Range("K132").Value = 7046.872499999999
$ws.Range("M132").Value = -4516.872499999999
